$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.031673198363605
$ws.Cells.Item(2, 4).Value = 1.051972450732113
$ws.Cells.Item(2, 5).Value = 1.031216568608252
$ws.Cells.Item(2, 6).Value = 1.057942626363359
$ws.Cells.Item(2, 9).Value = 1.044456028989515
$ws.Cells.Item(2, 10).Value = 1.036807629340518
$ws.Cells.Item(2, 11).Value = 1.054722472199143
$ws.Cells.Item(2, 12).Value = 1.034025327195363
$ws.Cells.Item(2, 13).Value = 1.060676215283342
$ws.Cells.Item(2, 14).Value = 1.038280015097679

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.032550661841594
$ws.Cells.Item(3, 4).Value = 1.052600584086589
$ws.Cells.Item(3, 5).Value = 1.031960201960181
$ws.Cells.Item(3, 6).Value = 1.05872414945783
$ws.Cells.Item(3, 9).Value = 1.044684897724369
$ws.Cells.Item(3, 10).Value = 1.037327433069893
$ws.Cells.Item(3, 11).Value = 1.055163832872516
$ws.Cells.Item(3, 12).Value = 1.034577652105992
$ws.Cells.Item(3, 13).Value = 1.061271756243174
$ws.Cells.Item(3, 14).Value = 1.038800557007971

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.033119011732486
$ws.Cells.Item(4, 4).Value = 1.053007335807203
$ws.Cells.Item(4, 5).Value = 1.032442256605858
$ws.Cells.Item(4, 6).Value = 1.059230538852148
$ws.Cells.Item(4, 9).Value = 1.044831912589931
$ws.Cells.Item(4, 10).Value = 1.037663713345012
$ws.Cells.Item(4, 11).Value = 1.055449038108559
$ws.Cells.Item(4, 12).Value = 1.034935251892772
$ws.Cells.Item(4, 13).Value = 1.061657141406283
$ws.Cells.Item(4, 14).Value = 1.039137314839653

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.033358081627758
$ws.Cells.Item(5, 4).Value = 1.053178405999587
$ws.Cells.Item(5, 5).Value = 1.032645119976384
$ws.Cells.Item(5, 6).Value = 1.059443588583592
$ws.Cells.Item(5, 9).Value = 1.044893458607108
$ws.Cells.Item(5, 10).Value = 1.037805068224434
$ws.Cells.Item(5, 11).Value = 1.055568844952112
$ws.Cells.Item(5, 12).Value = 1.035085635511789
$ws.Cells.Item(5, 13).Value = 1.061819162978619
$ws.Cells.Item(5, 14).Value = 1.039278870459218

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.033398230425651
$ws.Cells.Item(6, 4).Value = 1.053207133592862
$ws.Cells.Item(6, 5).Value = 1.032679193753868
$ws.Cells.Item(6, 6).Value = 1.059479370109206
$ws.Cells.Item(6, 9).Value = 1.044903777240175
$ws.Cells.Item(6, 10).Value = 1.037828801278795
$ws.Cells.Item(6, 11).Value = 1.05558895552964
$ws.Cells.Item(6, 12).Value = 1.035110888399894
$ws.Cells.Item(6, 13).Value = 1.061846367397906
$ws.Cells.Item(6, 14).Value = 1.039302637217238

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.03312220566508
$ws.Cells.Item(7, 4).Value = 1.053009621375792
$ws.Cells.Item(7, 5).Value = 1.032444966462532
$ws.Cells.Item(7, 6).Value = 1.059233384992603
$ws.Cells.Item(7, 9).Value = 1.044832735989689
$ws.Cells.Item(7, 10).Value = 1.037665602205488
$ws.Cells.Item(7, 11).Value = 1.055450639342318
$ws.Cells.Item(7, 12).Value = 1.034937261136422
$ws.Cells.Item(7, 13).Value = 1.061659306326292
$ws.Cells.Item(7, 14).Value = 1.039139206382528

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.031969621891589
$ws.Cells.Item(8, 4).Value = 1.052184666582431
$ws.Cells.Item(8, 5).Value = 1.031467700978415
$ws.Cells.Item(8, 6).Value = 1.058206601353139
$ws.Cells.Item(8, 9).Value = 1.044533599060467
$ws.Cells.Item(8, 10).Value = 1.036983312704818
$ws.Cells.Item(8, 11).Value = 1.054871710931107
$ws.Cells.Item(8, 12).Value = 1.034211944016171
$ws.Cells.Item(8, 13).Value = 1.060877473887378
$ws.Cells.Item(8, 14).Value = 1.038455947952507

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.029943068328784
$ws.Cells.Item(9, 4).Value = 1.050733426345546
$ws.Cells.Item(9, 5).Value = 1.029752403174543
$ws.Cells.Item(9, 6).Value = 1.056402661646181
$ws.Cells.Item(9, 9).Value = 1.043998258390907
$ws.Cells.Item(9, 10).Value = 1.035780567062035
$ws.Cells.Item(9, 11).Value = 1.05384867916397
$ws.Cells.Item(9, 12).Value = 1.032935497021225
$ws.Cells.Item(9, 13).Value = 1.059500088882167
$ws.Cells.Item(9, 14).Value = 1.037251494272955

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.028595106194215
$ws.Cells.Item(10, 4).Value = 1.049767679478493
$ws.Cells.Item(10, 5).Value = 1.028613521916232
$ws.Cells.Item(10, 6).Value = 1.055203772134233
$ws.Cells.Item(10, 9).Value = 1.043635885470655
$ws.Cells.Item(10, 10).Value = 1.034978495065795
$ws.Cells.Item(10, 11).Value = 1.053164798333805
$ws.Cells.Item(10, 12).Value = 1.032085720771503
$ws.Cells.Item(10, 13).Value = 1.058582132401929
$ws.Cells.Item(10, 14).Value = 1.036448283242481

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.028012170274084
$ws.Cells.Item(11, 4).Value = 1.049349937740171
$ws.Cells.Item(11, 5).Value = 1.028121497813149
$ws.Cells.Item(11, 6).Value = 1.054685550755346
$ws.Cells.Item(11, 9).Value = 1.043477684838876
$ws.Cells.Item(11, 10).Value = 1.034631146562221
$ws.Cells.Item(11, 11).Value = 1.052868246298284
$ws.Cells.Item(11, 12).Value = 1.031718056140097
$ws.Cells.Item(11, 13).Value = 1.058184738547989
$ws.Cells.Item(11, 14).Value = 1.036100441464191

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.027795754572539
$ws.Cells.Item(12, 4).Value = 1.049194836779281
$ws.Cells.Item(12, 5).Value = 1.027938907966818
$ws.Cells.Item(12, 6).Value = 1.054493198324519
$ws.Cells.Item(12, 9).Value = 1.043418728976826
$ws.Cells.Item(12, 10).Value = 1.034502120127232
$ws.Cells.Item(12, 11).Value = 1.05275803085997
$ws.Cells.Item(12, 12).Value = 1.031581534591941
$ws.Cells.Item(12, 13).Value = 1.058037143528246
$ws.Cells.Item(12, 14).Value = 1.035971231796864

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.027842171354942
$ws.Cells.Item(13, 4).Value = 1.04922810340372
$ws.Cells.Item(13, 5).Value = 1.02797806639706
$ws.Cells.Item(13, 6).Value = 1.054534452294264
$ws.Cells.Item(13, 9).Value = 1.04343138393454
$ws.Cells.Item(13, 10).Value = 1.03452979697608
$ws.Cells.Item(13, 11).Value = 1.05278167526598
$ws.Cells.Item(13, 12).Value = 1.03161081684943
$ws.Cells.Item(13, 13).Value = 1.058068802467362
$ws.Cells.Item(13, 14).Value = 1.035998947950012

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.027994278978509
$ws.Cells.Item(14, 4).Value = 1.049337115670438
$ws.Cells.Item(14, 5).Value = 1.028106401405556
$ws.Cells.Item(14, 6).Value = 1.054669648015445
$ws.Cells.Item(14, 9).Value = 1.043472815463438
$ws.Cells.Item(14, 10).Value = 1.034620481307015
$ws.Cells.Item(14, 11).Value = 1.05285913712917
$ws.Cells.Item(14, 12).Value = 1.031706770291427
$ws.Cells.Item(14, 13).Value = 1.058172537987284
$ws.Cells.Item(14, 14).Value = 1.036089761063099

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.02808801247505
$ws.Cells.Item(15, 4).Value = 1.049404290649123
$ws.Cells.Item(15, 5).Value = 1.028185495392453
$ws.Cells.Item(15, 6).Value = 1.054752964915567
$ws.Cells.Item(15, 9).Value = 1.043498317232247
$ws.Cells.Item(15, 10).Value = 1.034676354191658
$ws.Cells.Item(15, 11).Value = 1.05290685565765
$ws.Cells.Item(15, 12).Value = 1.031765896436967
$ws.Cells.Item(15, 13).Value = 1.058236454868173
$ws.Cells.Item(15, 14).Value = 1.036145713293647

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.028633809453156
$ws.Cells.Item(16, 4).Value = 1.049795412904391
$ws.Cells.Item(16, 5).Value = 1.02864619965632
$ws.Cells.Item(16, 6).Value = 1.055238184032499
$ws.Cells.Item(16, 9).Value = 1.043646357600533
$ws.Cells.Item(16, 10).Value = 1.03500154657291
$ws.Cells.Item(16, 11).Value = 1.053184470633038
$ws.Cells.Item(16, 12).Value = 1.032110127767856
$ws.Cells.Item(16, 13).Value = 1.058608508102693
$ws.Cells.Item(16, 14).Value = 1.03647136748538

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.028976372732414
$ws.Cells.Item(17, 4).Value = 1.050040870812146
$ws.Cells.Item(17, 5).Value = 1.028935488054631
$ws.Cells.Item(17, 6).Value = 1.055542793052339
$ws.Cells.Item(17, 9).Value = 1.043738874415825
$ws.Cells.Item(17, 10).Value = 1.035205519670804
$ws.Cells.Item(17, 11).Value = 1.053358497781615
$ws.Cells.Item(17, 12).Value = 1.032326134629367
$ws.Cells.Item(17, 13).Value = 1.058841911804735
$ws.Cells.Item(17, 14).Value = 1.036675630248469

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.029176255459183
$ws.Cells.Item(18, 4).Value = 1.050184083913115
$ws.Cells.Item(18, 5).Value = 1.029104332936328
$ws.Cells.Item(18, 6).Value = 1.055720553592569
$ws.Cells.Item(18, 9).Value = 1.043792713313705
$ws.Cells.Item(18, 10).Value = 1.035324489206766
$ws.Cells.Item(18, 11).Value = 1.053459963570618
$ws.Cells.Item(18, 12).Value = 1.032452156066825
$ws.Cells.Item(18, 13).Value = 1.058978060555182
$ws.Cells.Item(18, 14).Value = 1.036794768734819

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.029244422365507
$ws.Cells.Item(19, 4).Value = 1.050232922888463
$ws.Cells.Item(19, 5).Value = 1.029161922919813
$ws.Cells.Item(19, 6).Value = 1.055781180081912
$ws.Cells.Item(19, 9).Value = 1.043811049828135
$ws.Cells.Item(19, 10).Value = 1.035365053941292
$ws.Cells.Item(19, 11).Value = 1.053494553735485
$ws.Cells.Item(19, 12).Value = 1.032495130879722
$ws.Cells.Item(19, 13).Value = 1.059024485122374
$ws.Cells.Item(19, 14).Value = 1.03683539107592

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.028939611543211
$ws.Cells.Item(20, 4).Value = 1.05001453117949
$ws.Cells.Item(20, 5).Value = 1.028904438979268
$ws.Cells.Item(20, 6).Value = 1.055510102347974
$ws.Cells.Item(20, 9).Value = 1.043728961118666
$ws.Cells.Item(20, 10).Value = 1.035183635746886
$ws.Cells.Item(20, 11).Value = 1.053339830577822
$ws.Cells.Item(20, 12).Value = 1.032302956202001
$ws.Cells.Item(20, 13).Value = 1.058816868924103
$ws.Cells.Item(20, 14).Value = 1.036653715246869

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.027949483953549
$ws.Cells.Item(21, 4).Value = 1.049305012412667
$ws.Cells.Item(21, 5).Value = 1.028068605242937
$ws.Cells.Item(21, 6).Value = 1.054629832418009
$ws.Cells.Item(21, 9).Value = 1.043460620234847
$ws.Cells.Item(21, 10).Value = 1.034593777180629
$ws.Cells.Item(21, 11).Value = 1.052836328264033
$ws.Cells.Item(21, 12).Value = 1.031678513128951
$ws.Cells.Item(21, 13).Value = 1.058141990041935
$ws.Cells.Item(21, 14).Value = 1.03606301901379

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.027327603915003
$ws.Cells.Item(22, 4).Value = 1.048859297905099
$ws.Cells.Item(22, 5).Value = 1.027544066722595
$ws.Cells.Item(22, 6).Value = 1.054077171550681
$ws.Cells.Item(22, 9).Value = 1.0432907868398
$ws.Cells.Item(22, 10).Value = 1.034222877092415
$ws.Cells.Item(22, 11).Value = 1.052519393658916
$ws.Cells.Item(22, 12).Value = 1.031286164672822
$ws.Cells.Item(22, 13).Value = 1.057717753055044
$ws.Cells.Item(22, 14).Value = 1.035691592204907

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.027657211977234
$ws.Cells.Item(23, 4).Value = 1.049095542190076
$ws.Cells.Item(23, 5).Value = 1.027822040743066
$ws.Cells.Item(23, 6).Value = 1.054370071058728
$ws.Cells.Item(23, 9).Value = 1.043380924324611
$ws.Cells.Item(23, 10).Value = 1.034419500941978
$ws.Cells.Item(23, 11).Value = 1.052687440616581
$ws.Cells.Item(23, 12).Value = 1.031494130590271
$ws.Cells.Item(23, 13).Value = 1.057942640399596
$ws.Cells.Item(23, 14).Value = 1.035888495282891

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.028956222121896
$ws.Cells.Item(24, 4).Value = 1.050026432797238
$ws.Cells.Item(24, 5).Value = 1.028918468386894
$ws.Cells.Item(24, 6).Value = 1.055524873600777
$ws.Cells.Item(24, 9).Value = 1.043733440896011
$ws.Cells.Item(24, 10).Value = 1.035193524164226
$ws.Cells.Item(24, 11).Value = 1.05334826561148
$ws.Cells.Item(24, 12).Value = 1.032313429448366
$ws.Cells.Item(24, 13).Value = 1.058828184698395
$ws.Cells.Item(24, 14).Value = 1.036663617706897

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.03046644543571
$ws.Cells.Item(25, 4).Value = 1.051108306854579
$ws.Cells.Item(25, 5).Value = 1.030195036359506
$ws.Cells.Item(25, 6).Value = 1.056868373185893
$ws.Cells.Item(25, 9).Value = 1.044137626161749
$ws.Cells.Item(25, 10).Value = 1.036091553392187
$ws.Cells.Item(25, 11).Value = 1.054113491076681
$ws.Cells.Item(25, 12).Value = 1.033265285267232
$ws.Cells.Item(25, 13).Value = 1.059856130036993
$ws.Cells.Item(25, 14).Value = 1.037562922239367
